# "Generate Report for Handback"
#
# Updates the localization-status workbook to reflect a completed handback:
#   - Status text flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet (zh-cn / de-de columns) and on the per-language
#     "Status" column of the zh-cn / de-de sheets.
#   - The per-language sheets gain "Latest Target File" (a hyperlinked a.md),
#     "Latest Handback File" (the generated xliff name) and an updated
#     "Latest Handback DateTime" for both data rows.
#   - A couple of columns are widened so the longer status text / file names fit.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$aMdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71644b2cdc10670041e8d5993ba05b4b8506e79c/e2e/a.md"

# Helper: set a column's ColumnWidth so the saved OOXML `width` lands as close
# as possible to $targetStoredWidth (this host quantises ColumnWidth on 1/6
# character-width steps before persisting, so back-solve for the step whose
# stored value is nearest the desired width).
function Set-ColumnStoredWidth($ws, $col, $targetStoredWidth) {
    $n = [math]::Round(($targetStoredWidth - (5.0 / 6.0)) * 6.0)
    $ws.Columns.Item($col).ColumnWidth = $n / 6.0
}

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn (E) and de-de (F) status cells for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

Set-ColumnStoredWidth $wsOverview 5 29.9777047293527
Set-ColumnStoredWidth $wsOverview 6 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-25 18:38:42"
$wsZh.Range("K3").Value = "2016-08-25 18:38:42"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md") | Out-Null

Set-ColumnStoredWidth $wsZh 3 29.9777047293527
Set-ColumnStoredWidth $wsZh 10 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-25 18:38:49"
$wsDe.Range("K3").Value = "2016-08-25 18:38:49"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md") | Out-Null

Set-ColumnStoredWidth $wsDe 3 29.9777047293527
Set-ColumnStoredWidth $wsDe 10 40
